# Added thumbnails (journal cover) for Rural21 missed ones
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Fill in the new "Thumbnail" column (S) for every data row (2-62) with the
# journal-cover feed path that corresponds to that row's volume/issue.
$ws.Range("S2:S12").Value2  = "private://feeds/Rural21_2009-04_en.jpg"
$ws.Range("S13:S23").Value2 = "private://feeds/Rural21_2008-01_en.jpg"
$ws.Range("S24:S42").Value2 = "private://feeds/Rural21_2007-01_en.jpg"
$ws.Range("S43:S62").Value2 = "private://feeds/Rural21_2006-02_en.jpg"

# Widen column S (Thumbnail) now that it holds real data.
$ws.Columns.Item(19).ColumnWidth = 30.764322916666668

# Reflect the view state captured when this data was added: zoomed out a
# bit, scrolled so the frozen pane starts near the bottom of the data and
# the newly-filled column S is selected.
$av = $excel.ActiveWindow
$av.Zoom = 70

$av.FreezePanes = $false
$selTop = $ws.Range("A2").Select()
$av.FreezePanes = $true

$selFinal = $ws.Range("S43:S62").Select()
